$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 101
$ws.Cells.Item(101, 1).Value = 1
$ws.Cells.Item(101, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(101, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(101, 4).Value = 44418
$ws.Cells.Item(101, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(101, 5).Value = 15
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108005
$ws.Cells.Item(101, 10).Value = "Piña"
$ws.Cells.Item(101, 11).Value = "Caramelo"
$ws.Cells.Item(101, 12).Value = "Especial"
$ws.Cells.Item(101, 13).Value = 200
$ws.Cells.Item(101, 14).Value = 15000
$ws.Cells.Item(101, 15).Value = 16000
$ws.Cells.Item(101, 16).Value = 15500
$ws.Cells.Item(101, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(101, 18).Value = "Ecuador"
$ws.Cells.Item(101, 19).Value = 1550
$ws.Cells.Item(101, 20).Value = 10

# Row 102
$ws.Cells.Item(102, 1).Value = 1
$ws.Cells.Item(102, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(102, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(102, 4).Value = 44418
$ws.Cells.Item(102, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(102, 5).Value = 15
$ws.Cells.Item(102, 6).Value = "Fruta"
$ws.Cells.Item(102, 7).Value = 100108
$ws.Cells.Item(102, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(102, 9).Value = 100108005
$ws.Cells.Item(102, 10).Value = "Piña"
$ws.Cells.Item(102, 11).Value = "Caramelo"
$ws.Cells.Item(102, 12).Value = "Primera"
$ws.Cells.Item(102, 13).Value = 250
$ws.Cells.Item(102, 14).Value = 15000
$ws.Cells.Item(102, 15).Value = 16000
$ws.Cells.Item(102, 16).Value = 15500
$ws.Cells.Item(102, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(102, 18).Value = "Ecuador"
$ws.Cells.Item(102, 19).Value = 1292
$ws.Cells.Item(102, 20).Value = 12

# Row 103
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(103, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(103, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(103, 4).Value = 44418
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 5).Value = 15
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100108
$ws.Cells.Item(103, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(103, 9).Value = 100108005
$ws.Cells.Item(103, 10).Value = "Piña"
$ws.Cells.Item(103, 11).Value = "Caramelo"
$ws.Cells.Item(103, 12).Value = "Segunda"
$ws.Cells.Item(103, 13).Value = 270
$ws.Cells.Item(103, 14).Value = 15000
$ws.Cells.Item(103, 15).Value = 16000
$ws.Cells.Item(103, 16).Value = 15500
$ws.Cells.Item(103, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(103, 18).Value = "Ecuador"
$ws.Cells.Item(103, 19).Value = 1107
$ws.Cells.Item(103, 20).Value = 14

# Row 104
$ws.Cells.Item(104, 1).Value = 1
$ws.Cells.Item(104, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(104, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(104, 4).Value = 44418
$ws.Cells.Item(104, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(104, 5).Value = 15
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100108
$ws.Cells.Item(104, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(104, 9).Value = 100108005
$ws.Cells.Item(104, 10).Value = "Piña"
$ws.Cells.Item(104, 11).Value = "Caramelo"
$ws.Cells.Item(104, 12).Value = "Tercera"
$ws.Cells.Item(104, 13).Value = 270
$ws.Cells.Item(104, 14).Value = 15000
$ws.Cells.Item(104, 15).Value = 16000
$ws.Cells.Item(104, 16).Value = 15500
$ws.Cells.Item(104, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(104, 18).Value = "Ecuador"
$ws.Cells.Item(104, 19).Value = 969
$ws.Cells.Item(104, 20).Value = 16
